# Update countries & provincias Spain
#
# Applies the daily COVID-19 data refresh to the "Pais" worksheet:
#  - a few countries swap their ranking/position in the list
#    (Panama/Rumania, Uruguay/Guyana, Islas Malvinas/Montserrat)
#  - the numeric statistics (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) are refreshed
#    for the affected countries
#  - the "Datos actualizados a ..." timestamp footer is updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (row 1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 01:07"

# --- Helper to write a full data row (A:H) ---------------------------
# NOTE: positional parameters only (named parameter binding is not
# reliable in this PowerShell-style runtime).
function Set-CountryRow($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# row 4 - Estados Unidos
Set-CountryRow 4 "Estados Unidos" 6583888 34413 3873357 2514420 0 873 196111

# row 6 - Brasil
Set-CountryRow 6 "Brasil" 4238446 39114 3497337 611587 0 869 129522

# row 9 - Colombia
Set-CountryRow 9 "Colombia" 694664 7808 569479 102910 0 222 22275

# row 29 - Canada
Set-CountryRow 29 "Canada" 134785 491 118900 6722 0 8 9163

# row 35 - Egipto
Set-CountryRow 35 "Egipto" 100557 154 81597 13370 0 13 5590

# row 36 - Panama (was Rumania; Panama now ranks above Rumania)
Set-CountryRow 36 "Panama" 99715 673 72203 25385 0 11 2127

# row 37 - Rumania (was Panama)
Set-CountryRow 37 "Rumania" 99684 1380 41010 54609 0 47 4065

# row 48 - Japon
Set-CountryRow 48 "Japon" 73221 495 64835 6980 0 13 1406

# row 56 - Nigeria
Set-CountryRow 56 "Nigeria" 55829 197 43810 10944 0 5 1075

# row 84 - Bulgaria
Set-CountryRow 84 "Bulgaria" 17598 163 12619 4273 0 4 706

# row 92 - Noruega
Set-CountryRow 92 "Noruega" 11867 121 9348 2254 0 1 265

# row 108 - Montenegro
Set-CountryRow 108 "Montenegro" 6222 128 4393 1715 0 0 114

# row 139 - Trinidad yTobago
Set-CountryRow 139 "Trinidad yTobago" 2698 110 755 1900 0 4 43

# row 154 - Uruguay (was Guyana; Uruguay now ranks above Guyana)
Set-CountryRow 154 "Uruguay" 1759 18 1484 230 0 0 45

# row 155 - Guyana (was Uruguay)
Set-CountryRow 155 "Guyana" 1750 47 1088 613 0 1 49

# row 158 - Burkina Faso
Set-CountryRow 158 "Burkina Faso" 1486 10 1123 307 0 0 56

# row 165 - Vietnam
Set-CountryRow 165 "Vietnam" 1059 0 893 131 0 0 35

# row 214 - Islas Malvinas (was Montserrat; Islas Malvinas now ranks above Montserrat)
Set-CountryRow 214 "Islas Malvinas" 13 0 13 0 0 0 0

# row 215 - Montserrat (was Islas Malvinas)
Set-CountryRow 215 "Montserrat" 13 0 12 0 0 0 1
